$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.77588366666667
$ws.Range("H2").Value = 71.327651
$ws.Range("I2").Value = 0.201093431146956
$ws.Range("J2").Value = 0.2010934311469559
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 1675.318093630106
$ws.Range("R2").Value = 15077.86284267095
$ws.Range("S2").Value = 0.106111984283002
$ws.Range("T2").Value = 0.1061119842830019
$ws.Range("G3").Value = 23.77588366666667
$ws.Range("H3").Value = 71.327651
$ws.Range("I3").Value = 0.201093431146956
$ws.Range("J3").Value = 0.2010934311469559
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 233.7103663741468
$ws.Range("R3").Value = 2103.393297367321
$ws.Range("S3").Value = 0.01480284300501537
$ws.Range("T3").Value = 0.01480284300501537
$ws.Range("G4").Value = 23.77588366666667
$ws.Range("H4").Value = 71.327651
$ws.Range("I4").Value = 0.201093431146956
$ws.Range("J4").Value = 0.2010934311469559
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 1010.605775641733
$ws.Range("R4").Value = 9095.451980775593
$ws.Range("S4").Value = 0.06401016295887003
$ws.Range("T4").Value = 0.06401016295887003
$ws.Range("G5").Value = 23.77588366666667
$ws.Range("H5").Value = 71.327651
$ws.Range("I5").Value = 0.201093431146956
$ws.Range("J5").Value = 0.2010934311469559
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 255.2707101719242
$ws.Range("R5").Value = 2297.436391547318
$ws.Range("S5").Value = 0.01616844090006861
$ws.Range("T5").Value = 0.01616844090006861
$ws.Range("I6").Value = 0.2703947904457373
$ws.Range("J6").Value = 0.2703947904457373
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 2252.67072262555
$ws.Range("R6").Value = 20274.03650362995
$ws.Range("S6").Value = 0.1426805818088405
$ws.Range("T6").Value = 0.1426805818088405
$ws.Range("I7").Value = 0.2703947904457373
$ws.Range("J7").Value = 0.2703947904457373
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("S7").Value = 0.01990423858956006
$ws.Range("T7").Value = 0.01990423858956006
$ws.Range("I8").Value = 0.2703947904457373
$ws.Range("J8").Value = 0.2703947904457373
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 1358.883457153815
$ws.Range("R8").Value = 12229.95111438433
$ws.Range("S8").Value = 0.08606951754188695
$ws.Range("T8").Value = 0.08606951754188695
$ws.Range("I9").Value = 0.2703947904457373
$ws.Range("J9").Value = 0.2703947904457373
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 343.2427891363116
$ws.Range("R9").Value = 3089.185102226804
$ws.Range("S9").Value = 0.02174045250544982
$ws.Range("T9").Value = 0.02174045250544982
$ws.Range("G10").Value = 14.51831366666667
$ws.Range("H10").Value = 43.554941
$ws.Range("I10").Value = 0.1227940694288843
$ws.Range("J10").Value = 0.1227940694288843
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 1023.002716355985
$ws.Range("R10").Value = 9207.024447203868
$ws.Range("S10").Value = 0.06479536547248804
$ws.Range("T10").Value = 0.06479536547248804
$ws.Range("G11").Value = 14.51831366666667
$ws.Range("H11").Value = 43.554941
$ws.Range("I11").Value = 0.1227940694288843
$ws.Range("J11").Value = 0.1227940694288843
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("Q11").Value = 142.7110114493234
$ws.Range("R11").Value = 1284.399103043911
$ws.Range("S11").Value = 0.009039088553690169
$ws.Range("T11").Value = 0.009039088553690169
$ws.Range("G12").Value = 14.51831366666667
$ws.Range("H12").Value = 43.554941
$ws.Range("I12").Value = 0.1227940694288843
$ws.Range("J12").Value = 0.1227940694288843
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 617.1081525218726
$ws.Range("R12").Value = 5553.973372696853
$ws.Range("S12").Value = 0.03908664917444105
$ws.Range("T12").Value = 0.03908664917444105
$ws.Range("G13").Value = 14.51831366666667
$ws.Range("H13").Value = 43.554941
$ws.Range("I13").Value = 0.1227940694288843
$ws.Range("J13").Value = 0.1227940694288843
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 155.8764457358375
$ws.Range("R13").Value = 1402.888011622538
$ws.Range("S13").Value = 0.009872966228265041
$ws.Range("T13").Value = 0.009872966228265041
$ws.Range("G14").Value = 47.96923
$ws.Range("H14").Value = 143.90769
$ws.Range("I14").Value = 0.4057177089784224
$ws.Range("J14").Value = 0.4057177089784224
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 3380.051824074681
$ws.Range("R14").Value = 30420.46641667212
$ws.Range("S14").Value = 0.2140871082307634
$ws.Range("T14").Value = 0.2140871082307634
$ws.Range("G15").Value = 47.96923
$ws.Range("H15").Value = 143.90769
$ws.Range("I15").Value = 0.4057177089784224
$ws.Range("J15").Value = 0.4057177089784224
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 471.5242754027767
$ws.Range("R15").Value = 4243.71847862499
$ws.Range("S15").Value = 0.02986559787710407
$ws.Range("T15").Value = 0.02986559787710407
$ws.Range("G16").Value = 47.96923
$ws.Range("H16").Value = 143.90769
$ws.Range("I16").Value = 0.4057177089784224
$ws.Range("J16").Value = 0.4057177089784224
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 2038.95600982654
$ws.Range("R16").Value = 18350.60408843886
$ws.Range("S16").Value = 0.1291442317080448
$ws.Range("T16").Value = 0.1291442317080448
$ws.Range("G17").Value = 47.96923
$ws.Range("H17").Value = 143.90769
$ws.Range("I17").Value = 0.4057177089784224
$ws.Range("J17").Value = 0.4057177089784224
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 515.0235246847134
$ws.Range("R17").Value = 4635.21172216242
$ws.Range("S17").Value = 0.03262077116251023
$ws.Range("T17").Value = 0.03262077116251023
